$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the article title for student 212241811426 (row 4) in column D
# from "Cyberpunk and the Dilemmas of Postmodern Narrative: The Example of William Gibson"
# to "Claire Sponsler_1992_Contemporary Literature_Cyberpunk and the Dilemmas of Postmodern Narrative: The Example of William Gibson"
$ws.Range("D4").Value = "Claire Sponsler_1992_Contemporary Literature_Cyberpunk and the Dilemmas of Postmodern Narrative: The Example of William Gibson"

# Reflect the user's final selection/active cell being E5
$ws.Range("E5").Select()

$wb.Save()
